$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 20
$ws.Range("H20").Value = 1460.5
$ws.Range("I20").Value = 1460.5
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1460.5
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -1230.5

# Row 35
$ws.Range("H35").Value = 1460.5
$ws.Range("I35").Value = 1460.5
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1460.5
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -1081.5

# Row 74
$ws.Range("H74").Value = 3443.75
$ws.Range("I74").Value = 3675
$ws.Range("J74").Value = 2750
$ws.Range("K74").Value = 3675
$ws.Range("L74").Value = 2750
$ws.Range("M74").Value = -2739
$ws.Range("N74").Value = -4622

# Row 77
$ws.Range("H77").Value = 3443.75
$ws.Range("I77").Value = 3675
$ws.Range("J77").Value = 2750
$ws.Range("K77").Value = 18375
$ws.Range("L77").Value = 13750
$ws.Range("M77").Value = -13695
$ws.Range("N77").Value = -23110

# Row 100
$ws.Range("H100").Value = 1914.2858
$ws.Range("I100").Value = 1500
$ws.Range("J100").Value = 2950
$ws.Range("K100").Value = 1500
$ws.Range("L100").Value = 2950
$ws.Range("M100").Value = -959
$ws.Range("N100").Value = -4032

# Row 141
$ws.Range("H141").Value = 1957.0333
$ws.Range("I141").Value = 1020.04
$ws.Range("J141").Value = 6642
$ws.Range("K141").Value = 3060.12
$ws.Range("L141").Value = 19926
$ws.Range("M141").Value = 2119.88
$ws.Range("N141").Value = -30286


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 9909.241
$ws.Range("I32").Value = 9242.544
$ws.Range("J32").Value = 12931.6
$ws.Range("K32").Value = 9242.544
$ws.Range("L32").Value = 12931.6
$ws.Range("M32").Value = -8955.544
$ws.Range("N32").Value = -13505.6

# Row 37
$ws.Range("H37").Value = 47816.668
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 47816.668
$ws.Range("K37").Value = 0
$ws.Range("L37").ClearContents()
$ws.Range("M37").Value = 47816.668
$ws.Range("N37").Value = -48362.668

# Row 61
$ws.Range("H61").Value = 1252.8555
$ws.Range("I61").Value = 1104.1724
$ws.Range("K61").Value = 1104.1724
$ws.Range("M61").Value = -892.1723999999999

# Row 74
$ws.Range("H74").Value = 1698.3846
$ws.Range("I74").Value = 1560.2236
$ws.Range("K74").Value = 1560.2236
$ws.Range("M74").Value = -686.2236

# Row 77
$ws.Range("H77").Value = 1698.3846
$ws.Range("I77").Value = 1560.2236
$ws.Range("K77").Value = 7801.118
$ws.Range("M77").Value = -3433.118

# Row 132
$ws.Range("H132").Value = 7354429
$ws.Range("I132").Value = 10870514
$ws.Range("J132").Value = 2614.318
$ws.Range("K132").Value = 32611542
$ws.Range("L132").Value = 7842.954000000001
$ws.Range("M132").Value = -32609012
$ws.Range("N132").Value = -12902.954

# Row 136
$ws.Range("H136").Value = 1252.8555
$ws.Range("I136").Value = 1104.1724
$ws.Range("K136").Value = 3312.5172
$ws.Range("M136").Value = -762.5171999999998


$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 1033.3334
$ws.Range("I22").Value = 1375
$ws.Range("J22").Value = 350
$ws.Range("K22").Value = 1375
$ws.Range("L22").Value = 350
$ws.Range("M22").Value = -1202
$ws.Range("N22").Value = -696

# Row 35
$ws.Range("H35").Value = 47000
$ws.Range("J35").Value = 47000
$ws.Range("L35").Value = 47000
$ws.Range("N35").Value = -47620


$ws = $wb.Worksheets.Item("CRP")
# Row 41
$ws.Range("H41").Value = 31942.8
$ws.Range("J41").Value = 38303.5
$ws.Range("L41").Value = 38303.5
$ws.Range("N41").Value = -39159.5

# Row 47
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").ClearContents()
$ws.Range("N47").Value = 0

# Row 48
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").ClearContents()
$ws.Range("N48").Value = 0

# Row 50
$ws.Range("H50").Value = 39602.25
$ws.Range("J50").Value = 39602.25
$ws.Range("L50").Value = 39602.25
$ws.Range("N50").Value = -40852.25

# Row 51
$ws.Range("H51").Value = 100031390
$ws.Range("J51").Value = 39237.25
$ws.Range("L51").Value = 39237.25
$ws.Range("N51").Value = -40709.25

# Row 59
$ws.Range("H59").Value = 30902.777
$ws.Range("J59").Value = 29752.625
$ws.Range("L59").Value = 29752.625
$ws.Range("N59").Value = -32042.625

# Row 60
$ws.Range("H60").Value = 23531.191
$ws.Range("I60").Value = 23103
$ws.Range("J60").Value = 23552.6
$ws.Range("K60").Value = 23103
$ws.Range("L60").Value = 23552.6
$ws.Range("M60").Value = -22592
$ws.Range("N60").Value = -24574.6

# Row 61
$ws.Range("H61").Value = 100031390
$ws.Range("J61").Value = 39237.25
$ws.Range("L61").Value = 39237.25
$ws.Range("N61").Value = -39933.25

# Row 62
$ws.Range("H62").Value = 3165.5557
$ws.Range("I62").Value = 2696.8
$ws.Range("J62").Value = 3751.5
$ws.Range("K62").Value = 2696.8
$ws.Range("L62").Value = 3751.5
$ws.Range("M62").Value = -2072.8
$ws.Range("N62").Value = -4999.5

# Row 65
$ws.Range("H65").Value = 3165.5557
$ws.Range("I65").Value = 2696.8
$ws.Range("J65").Value = 3751.5
$ws.Range("K65").Value = 13484
$ws.Range("L65").Value = 18757.5
$ws.Range("M65").Value = -10364
$ws.Range("N65").Value = -24997.5

# Row 103
$ws.Range("H103").Value = 17555.555
$ws.Range("I103").Value = 9000
$ws.Range("K103").Value = 9000
$ws.Range("M103").Value = -7828

# Row 132
$ws.Range("H132").Value = 28726.98
$ws.Range("I132").Value = 1363.762
$ws.Range("J132").Value = 143652.5
$ws.Range("K132").Value = 4091.286
$ws.Range("L132").Value = 430957.5
$ws.Range("M132").Value = -1561.286
$ws.Range("N132").Value = -436017.5

# Row 134
$ws.Range("H134").Value = 265567.22
$ws.Range("I134").Value = 878.9737
$ws.Range("J134").Value = 936110.75
$ws.Range("K134").Value = 2636.9211
$ws.Range("L134").Value = 2808332.25
$ws.Range("M134").Value = -101.9211
$ws.Range("N134").Value = -2813402.25


$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 2222.5386
$ws.Range("I100").Value = 2149.125
$ws.Range("J100").Value = 2340
$ws.Range("K100").Value = 2149.125
$ws.Range("L100").Value = 2340
$ws.Range("M100").Value = -1608.125
$ws.Range("N100").Value = -3422

# Row 132
$ws.Range("H132").Value = 2220.0195
$ws.Range("I132").Value = 1412
$ws.Range("K132").Value = 4236
$ws.Range("M132").Value = -1706


$ws = $wb.Worksheets.Item("WVR")
# Row 17
$ws.Range("H17").Value = 2850
$ws.Range("I17").Value = 2850
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 2850
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -2678

# Row 132
$ws.Range("H132").Value = 1861.0164
$ws.Range("I132").Value = 1816.7906
$ws.Range("J132").Value = 1966.6666
$ws.Range("K132").Value = 5450.3718
$ws.Range("L132").Value = 5899.9998
$ws.Range("M132").Value = -2920.3718
$ws.Range("N132").Value = -10959.9998

